$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2,3,4 (rat_brainstem_* entries), which shifts the old row 5
# (scaffold_context_info.json / application/x.vnd... / {"version"...}) up to row 2.
$ws.Range("A2:C4").EntireRow.Delete() | Out-Null

# Update the selection to match the saved view state.
$ws.Range("A10").Select() | Out-Null
